# Refresh the Cd86-Ctla4 LR-pair sheet with new TPM-based NATMI output.
# The sending/target cluster pairing grows from a single combination
# (ECs -> Resolving-Mac) to the full 2x3 cross of
# Sending clusters {ECs, Resolving-Mac} x Target clusters {ECs, MuSCs, Resolving-Mac},
# with all of the downstream expression / specificity statistics recomputed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> ECs
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Cd86"
$ws.Cells.Item(2, 3).Value = "Ctla4"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 1
$ws.Cells.Item(2, 6).Value = 0.3333333333333333
$ws.Cells.Item(2, 7).Value = 0.02345233333333334
$ws.Cells.Item(2, 8).Value = 0.070357
$ws.Cells.Item(2, 9).Value = 0.0002537772683371841
$ws.Cells.Item(2, 10).Value = 0.0002537772683371841
$ws.Cells.Item(2, 11).Value = 1
$ws.Cells.Item(2, 12).Value = 0.3333333333333333
$ws.Cells.Item(2, 13).Value = 0.04641433333333334
$ws.Cells.Item(2, 14).Value = 0.139243
$ws.Cells.Item(2, 15).Value = 0.1185678497650663
$ws.Cells.Item(2, 16).Value = 0.1185678497650663
$ws.Cells.Item(2, 17).Value = 0.001088524416777778
$ws.Cells.Item(2, 18).Value = 0.009796719751000001
$ws.Cells.Item(2, 19).Value = 0.00003008982502599217
$ws.Cells.Item(2, 20).Value = 0.00003008982502599217

# Row 3: ECs -> MuSCs (new target cluster)
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Cd86"
$ws.Cells.Item(3, 3).Value = "Ctla4"
$ws.Cells.Item(3, 4).Value = "MuSCs"
$ws.Cells.Item(3, 5).Value = 1
$ws.Cells.Item(3, 6).Value = 0.3333333333333333
$ws.Cells.Item(3, 7).Value = 0.02345233333333334
$ws.Cells.Item(3, 8).Value = 0.070357
$ws.Cells.Item(3, 9).Value = 0.0002537772683371841
$ws.Cells.Item(3, 10).Value = 0.0002537772683371841
$ws.Cells.Item(3, 11).Value = 1
$ws.Cells.Item(3, 12).Value = 0.3333333333333333
$ws.Cells.Item(3, 13).Value = 0.06025633333333334
$ws.Cells.Item(3, 14).Value = 0.180769
$ws.Cells.Item(3, 15).Value = 0.1539279650264737
$ws.Cells.Item(3, 16).Value = 0.1539279650264737
$ws.Cells.Item(3, 17).Value = 0.001413151614777778
$ws.Cells.Item(3, 18).Value = 0.012718364533
$ws.Cells.Item(3, 19).Value = 0.0000390634184851201
$ws.Cells.Item(3, 20).Value = 0.0000390634184851201

# Row 4: ECs -> Resolving-Mac
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Cd86"
$ws.Cells.Item(4, 3).Value = "Ctla4"
$ws.Cells.Item(4, 4).Value = "Resolving-Mac"
$ws.Cells.Item(4, 5).Value = 1
$ws.Cells.Item(4, 6).Value = 0.3333333333333333
$ws.Cells.Item(4, 7).Value = 0.02345233333333334
$ws.Cells.Item(4, 8).Value = 0.070357
$ws.Cells.Item(4, 9).Value = 0.0002537772683371841
$ws.Cells.Item(4, 10).Value = 0.0002537772683371841
$ws.Cells.Item(4, 11).Value = 2
$ws.Cells.Item(4, 12).Value = 0.6666666666666666
$ws.Cells.Item(4, 13).Value = 0.2847873333333333
$ws.Cells.Item(4, 14).Value = 0.8543620000000001
$ws.Cells.Item(4, 15).Value = 0.7275041852084601
$ws.Cells.Item(4, 16).Value = 0.7275041852084601
$ws.Cells.Item(4, 17).Value = 0.006678927470444445
$ws.Cells.Item(4, 18).Value = 0.06011034723400001
$ws.Cells.Item(4, 19).Value = 0.0001846240248260718
$ws.Cells.Item(4, 20).Value = 0.0001846240248260718

# Row 5: Resolving-Mac -> ECs
$ws.Cells.Item(5, 1).Value = "Resolving-Mac"
$ws.Cells.Item(5, 2).Value = "Cd86"
$ws.Cells.Item(5, 3).Value = "Ctla4"
$ws.Cells.Item(5, 4).Value = "ECs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 92.38960533333334
$ws.Cells.Item(5, 8).Value = 277.168816
$ws.Cells.Item(5, 9).Value = 0.9997462227316628
$ws.Cells.Item(5, 10).Value = 0.9997462227316628
$ws.Cells.Item(5, 11).Value = 1
$ws.Cells.Item(5, 12).Value = 0.3333333333333333
$ws.Cells.Item(5, 13).Value = 0.04641433333333334
$ws.Cells.Item(5, 14).Value = 0.139243
$ws.Cells.Item(5, 15).Value = 0.1185678497650663
$ws.Cells.Item(5, 16).Value = 0.1185678497650663
$ws.Cells.Item(5, 17).Value = 4.288201938476445
$ws.Cells.Item(5, 18).Value = 38.593817446288
$ws.Cells.Item(5, 19).Value = 0.1185377599400403
$ws.Cells.Item(5, 20).Value = 0.1185377599400403

# Row 6: Resolving-Mac -> MuSCs (new target cluster)
$ws.Cells.Item(6, 1).Value = "Resolving-Mac"
$ws.Cells.Item(6, 2).Value = "Cd86"
$ws.Cells.Item(6, 3).Value = "Ctla4"
$ws.Cells.Item(6, 4).Value = "MuSCs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 92.38960533333334
$ws.Cells.Item(6, 8).Value = 277.168816
$ws.Cells.Item(6, 9).Value = 0.9997462227316628
$ws.Cells.Item(6, 10).Value = 0.9997462227316628
$ws.Cells.Item(6, 11).Value = 1
$ws.Cells.Item(6, 12).Value = 0.3333333333333333
$ws.Cells.Item(6, 13).Value = 0.06025633333333334
$ws.Cells.Item(6, 14).Value = 0.180769
$ws.Cells.Item(6, 15).Value = 0.1539279650264737
$ws.Cells.Item(6, 16).Value = 0.1539279650264737
$ws.Cells.Item(6, 17).Value = 5.567058855500445
$ws.Cells.Item(6, 18).Value = 50.103529699504
$ws.Cells.Item(6, 19).Value = 0.1538889016079886
$ws.Cells.Item(6, 20).Value = 0.1538889016079886

# Row 7: Resolving-Mac -> Resolving-Mac
$ws.Cells.Item(7, 1).Value = "Resolving-Mac"
$ws.Cells.Item(7, 2).Value = "Cd86"
$ws.Cells.Item(7, 3).Value = "Ctla4"
$ws.Cells.Item(7, 4).Value = "Resolving-Mac"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 92.38960533333334
$ws.Cells.Item(7, 8).Value = 277.168816
$ws.Cells.Item(7, 9).Value = 0.9997462227316628
$ws.Cells.Item(7, 10).Value = 0.9997462227316628
$ws.Cells.Item(7, 11).Value = 2
$ws.Cells.Item(7, 12).Value = 0.6666666666666666
$ws.Cells.Item(7, 13).Value = 0.2847873333333333
$ws.Cells.Item(7, 14).Value = 0.8543620000000001
$ws.Cells.Item(7, 15).Value = 0.7275041852084601
$ws.Cells.Item(7, 16).Value = 0.7275041852084601
$ws.Cells.Item(7, 17).Value = 26.31138933059911
$ws.Cells.Item(7, 18).Value = 236.802503975392
$ws.Cells.Item(7, 19).Value = 0.7273195611836339
$ws.Cells.Item(7, 20).Value = 0.7273195611836339
